# Apply updates to the "Initial placement strategy" comparison sheet (Tabelle3)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle3")

# --- Column A / B existing values get refreshed with new benchmark numbers ---
$ws.Range("B5").Value = "{0=192, 1=14, 2=6, 3=14, 4=0, 5=0, 6=0, 7=1}"
$ws.Range("B6").Value = "{0=0, 1=4, 2=6, 3=0, 4=195, 5=0, 6=0, 7=0, 8=22}"
$ws.Range("B7").Value = "{0=0, 1=16, 2=10, 3=12, 4=168, 5=5, 6=0, 7=0, 8=16, 9=0}"

# Row 9 label switches from "w/o G4 (<-0)" to the new strategy name
$ws.Range("A9").Value = "Best linkable player first"

$ws.Range("B10").Value = "{0=0, 1=14, 2=12, 3=12, 4=168, 5=5, 6=0, 7=0, 8=16, 9=0}"
$ws.Range("B11").Value = "{0=0, 1=15, 2=8, 3=3, 4=180, 5=5, 6=0, 7=0, 8=16, 9=0}"

# --- New column C (second strategy comparison column) ---
$ws.Range("C5").Value = "{0=166, 1=19, 2=28, 3=6, 4=4, 5=3, 6=0, 7=1}"
$ws.Range("C6").Value = "{0=0, 1=4, 2=6, 3=0, 4=195, 5=0, 6=0, 7=0, 8=22}"
$ws.Range("C7").Value = "{0=0, 1=15, 2=24, 3=18, 4=156, 5=0, 6=6, 7=0, 8=8, 9=0}"
$ws.Range("C9").Value = "Least linkable player first"
$ws.Range("C9").Font.Bold = $true

# --- New column D (third data column, next to C's Shift/Break labels) ---
$ws.Range("D10").Value = "{0=0, 1=16, 2=16, 3=21, 4=160, 5=0, 6=6, 7=0, 8=8, 9=0}"
$ws.Range("D11").Value = "{0=0, 1=16, 2=12, 3=9, 4=176, 5=0, 6=6, 7=0, 8=8, 9=0}"

# Widen column H for upcoming notes/data, and move the active selection
# (input is pre-compensated so the saved OOXML column width lands on exactly 20)
$ws.Columns.Item(8).ColumnWidth = 19.17
$ws.Range("J9").Select()
